$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.122.15"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "2.955.44"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'381.14"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").Value = "'102.32"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("D7").Value = "'0.538"
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "'36.54"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "3.422.49"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("D14").Value = "'18.06"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").Value = "'7.41"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("D16").Value = "2.943.98"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'0.984"
$ws.Range("E17").Value = "  +4.17%  "
$ws.Range("D18").Value = "51.094.36"
$ws.Range("E19").Value = "  -5.94%  "
$ws.Range("D20").Value = "'7.09"
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("D21").Value = "'12.56"
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("D22").Value = "0.0₃0953"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'68.48"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'261.63"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("D26").Value = "'8.38"
$ws.Range("E26").Value = "  +14.01%  "
$ws.Range("D27").Value = "'7.59"
$ws.Range("E27").Value = "  +3.19%  "
$ws.Range("D28").Value = "'0.168"
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").Value = "'4.11"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.112"
$ws.Range("E31").Value = "  +9.20%  "
$ws.Range("D32").Value = "'25.66"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "'9.81"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").Value = "'0.0458"
$ws.Range("E34").Value = "  +6.87%  "
$ws.Range("D35").Value = "'33.87"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'50.42"
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("B37").Value = "Toncoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D37").Value = "'2.05"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").Value = "'16.81"
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("D41").Value = "'2.54"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("D44").Value = "'121.57"
$ws.Range("E44").Value = "  -2.15%  "
$ws.Range("D45").Value = "'21.16"
$ws.Range("E45").Value = "  -2.84%  "
$ws.Range("D46").Value = "'2.06"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "'0.274"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("D49").Value = "2.006.20"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").Value = "'0.0337"
$ws.Range("E51").Value = "  +4.76%  "
